# Add source_project column + source_project list sheet (HCA lookup)
# per HuBMAP scrnaseq-hca metadata template change.

$wb = $excel.ActiveWorkbook
$tsv = $wb.Worksheets.Item("Export as TSV")

# ---- 1. Capture existing header comments (columns A..AO) before we move anything ----
$oldComments = @(
    'HuBMAP Display ID of the donor of the assayed tissue.',
    'HuBMAP Display ID of the assayed tissue.',
    'Start date and time of assay, typically a date-time stamped folder generated by the acquisition instrument. YYYY-MM-DD hh:mm, where YYYY is the year, MM is the month with leading 0s, and DD is the day with leading 0s, hh is the hour with leading zeros, mm are the minutes with leading zeros.',
    'DOI for protocols.io referring to the protocol for this assay.',
    'Name of the person responsible for executing the assay.',
    'Email address for the operator.',
    'Name of the principal investigator responsible for the data.',
    'Email address for the principal investigator.',
    'Each assay is placed into one of the following 3 general categories: generation of images of microscopic entities, identification & quantitation of molecules by mass spectrometry, and determination of nucleotide sequence.',
    'The specific type of assay being executed.',
    'Analytes are the target molecules being measured with the assay.',
    'Specifies whether or not a specific molecule(s) is/are targeted for detection/measurement by the assay. For example, an antibody targets a specific protein.',
    'An acquisition instrument is the device that contains the signal detection hardware and signal processing software. Assays generate signals such as light of various intensities or color or signals representing the molecular mass.',
    'Manufacturers of an acquisition instrument may offer various versions (models) of that instrument with different features or sensitivities. Differences in features or sensitivities may be relevant to processing or interpretation of the data.',
    'Link to a protocols document answering the question: How were single cells separated into a single-cell suspension?',
    'The type of single cell entity derived from isolation protocol',
    'The method by which tissues are dissociated into single cells in suspension.',
    'The method by which specific cell populations are sorted or enriched.',
    'A quality metric by visual inspection prior to cell lysis or defined by known parameters such as wells with several cells or no cells. This can be captured at a high level.',
    'Total number of cell/nuclei yielded post dissociation and enrichment',
    'Number of cell/nuclei input to the assay',
    'The kit used for the RNA sequencing assay',
    'A link to the protocol document containing the library construction method (including version) that was used, e.g. "Smart-Seq2", "Drop-Seq", "10X v3".',
    'Whether the library was generated for single-end or paired end sequencing',
    'Adapter sequence to be used for adapter trimming',
    'An id for the library. The id may be text and/or numbers',
    'Is the sequencing reaction run in replicate, TRUE or FALSE',
    'Which read file contains the cell barcode',
    'Position(s) in the read at which the cell barcode starts.',
    'Length of the cell barcode in base pairs',
    'Number of PCR cycles to amplify cDNA',
    'Number of PCR cycles performed for library indexing',
    'Total number of ng of library after final pcr amplification step. This is the concentration (ng/ul) * volume (ul)',
    'Units of final library yield',
    'Average size of sequencing library fragments estimated via gel electrophoresis or bioanalyzer/tapestation.',
    'Reagent kit used for sequencing',
    'Slash-delimited list of the number of sequencing cycles for, for example, Read1, i7 index, i5 index, and Read2.',
    'Percent of bases with Quality scores above Q30',
    'Percent PhiX loaded to the run',
    'Relative path to file with ORCID IDs for contributors for this dataset.',
    'Relative path to file or directory with instrument data. Downstream processing will depend on filename extension conventions.',
)

# ---- 2. Delete the existing comments so they do not linger at their old (pre-shift) cells ----
for ($i = 1; $i -le $oldComments.Length; $i++) {
    $cell = $tsv.Cells.Item(1, $i)
    if ($cell.Comment -ne $null) {
        $cell.Comment.Delete()
    }
}

# ---- 3. Insert the new first column (shifts data, validations; data validation sqref shifts automatically) ----
$tsv.Columns.Item(1).Insert()

# ---- 4. New header + value for the inserted column (copy header formatting from B1) ----
$tsv.Range("A1").Value = "source_project"
$tsv.Range("B1").Copy()
$tsv.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- 5. Re-create the header comments, shifted one column to the right (B..AP) ----
for ($i = 1; $i -le $oldComments.Length; $i++) {
    $cell = $tsv.Cells.Item(1, $i + 1)
    $cell.AddComment($oldComments[$i - 1])
}

# ---- 6. New comment describing source_project on the inserted column A ----
$tsv.Range("A1").AddComment('External source (outside of HuBMAP) of the project, eg. HCA (The Human Cell Atlas Consortium).')

# ---- 7. Add the "source_project list" lookup sheet right after "Export as TSV" ----
$sourceProjectList = $wb.Worksheets.Add($null, $tsv)
$sourceProjectList.Name = "source_project list"
$sourceProjectList.Range("A1").Value = "HCA"

# ---- 8. Data validation for the new column, restricting to the source_project list ----
$validationRange = $tsv.Range("A2:A1048576")
$validationRange.Validation.Add(3, 1, 1, "='source_project list'!`$A`$1:`$A`$1")
$validationRange.Validation.ErrorTitle = "Value must come from list"
$validationRange.Validation.ErrorMessage = "Value must be one of: HCA."
$validationRange.Validation.ShowInput = $true
$validationRange.Validation.ShowError = $true
